# Generate Report for Handback
# Adds "Latest Target File" hyperlink, "Latest Handback File", a new
# "Latest Handback DateTime" and an "Error Detail" message to row 8
# (the cc917249-... entry) on both the zh-cn and de-de sheets, because
# the handed-back file version is not the latest. Also widens the
# "Error Detail" column so the message is readable.

$wb = $excel.ActiveWorkbook

$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/67c43828abe37023df6ce03e9390d058e627c9d0/e2e/cc917249-a881-4222-93e5-4a6e275b265d.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0db03b7f0cad8a8f4196144d388be4e47a544c3a/e2e/cc917249-a881-4222-93e5-4a6e275b265d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/67c43828abe37023df6ce03e9390d058e627c9d0/e2e/cc917249-a881-4222-93e5-4a6e275b265d.md."

$sheetConfigs = @(
    @{ Name = "zh-cn"; TargetFile = "cc917249-a881-4222-93e5-4a6e275b265d.48d27acd5d3a61f8ae729ea7279ec469144f2a8b.zh-cn.xlf"; HandbackDate = "2016-08-17 14:43:29" },
    @{ Name = "de-de"; TargetFile = "cc917249-a881-4222-93e5-4a6e275b265d.48d27acd5d3a61f8ae729ea7279ec469144f2a8b.de-de.xlf"; HandbackDate = "2016-08-17 14:43:37" }
)

foreach ($cfg in $sheetConfigs) {
    $ws = $wb.Worksheets.Item($cfg.Name)

    # Latest Target File (I8) becomes a hyperlink to the latest handback md
    $ws.Range("I8").Value = "cc917249-a881-4222-93e5-4a6e275b265d.md"
    $ws.Hyperlinks.Add($ws.Range("I8"), $latestUrl, "", "", "cc917249-a881-4222-93e5-4a6e275b265d.md")

    # Latest Handback File (J8)
    $ws.Range("J8").Value = $cfg.TargetFile

    # Latest Handback DateTime (K8)
    $ws.Range("K8").Value = $cfg.HandbackDate

    # Error Detail (P8)
    $ws.Range("P8").Value = $errorDetail

    # Widen the Error Detail column so the message is readable
    $ws.Columns.Item(16).ColumnWidth = 40
}
